{"js": "// Replace the seven \"Heading 3\" job-title/company/date lines in the\n// PROFESSIONAL EXPERIENCE section of the resume.\n//\n// Each entry maps the OLD full paragraph text to the NEW full paragraph\n// text. We search the whole document body for an exact (case-sensitive)\n// match of the old text and replace the hit range's text, which preserves\n// the run/paragraph formatting (Heading3 style, bold, etc.) already on\n// that line.\nconst replacements = [\n  [\n    \"PARTNER & SENIOR DATA ARCHITECT - Siege Analytics, Washington, DC | January 2014 \\u2013 Present\",\n    \"PARTNER - Siege Analytics, Washington, DC | January 2014 \\u2013 Present\",\n  ],\n  [\n    \"PRINCIPAL TECHNICAL ARCHITECT - Clarity and Rigour, Washington, DC | 2012 \\u2013 2014\",\n    \"DATA PRODUCTS MANAGER - Helm/Murmuration, Washington, DC | 2012 \\u2013 2014\",\n  ],\n  [\n    \"DIRECTOR OF TECHNOLOGY - Helm, Washington, DC | 2010 \\u2013 2012\",\n    \"SOFTWARE ENGINEER - Mautinoa Technologies, Washington, DC | 2010 \\u2013 2012\",\n  ],\n  [\n    \"SENIOR TECHNICAL ANALYST - GSD&M, Austin, TX | 2008 \\u2013 2010\",\n    \"SENIOR ANALYST - Myers Research, Washington, DC | 2008 \\u2013 2010\",\n  ],\n  [\n    \"TECHNICAL COORDINATOR - Progressive Change Campaign Committee, Washington, DC | 2006 \\u2013 2008\",\n    \"RESEARCH DIRECTOR - Progressive Change Campaign Committee, Washington, DC | 2006 \\u2013 2008\",\n  ],\n  [\n    \"TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 \\u2013 2004\",\n    \"INTERIM TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 \\u2013 2004\",\n  ],\n  [\n    \"TECHNICAL COORDINATOR - The Feldman Group, Washington, DC | 2000 \\u2013 2001\",\n    \"FIELD DIRECTOR - The Feldman Group, Washington, DC | 2000 \\u2013 2001\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(\"Could not find text: \" + oldText);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the seven \"Heading 3\" job-title/company/date lines in the\n# PROFESSIONAL EXPERIENCE section of the resume.\n#\n# Each pair is the OLD full paragraph text and the NEW full paragraph\n# text. We run Find/Replace (wdReplaceAll) scoped to the exact, unique\n# full-line text (title + company + location + dates) so the two\n# \"TECHNICAL COORDINATOR - ...\" paragraphs -- which share a prefix but\n# differ in company/location -- are each retargeted to their own new\n# title without touching the other.\n\n$d = $word.ActiveDocument\n$en_dash = [char]8211\n\n$replacements = @(\n  @(\"PARTNER & SENIOR DATA ARCHITECT - Siege Analytics, Washington, DC | January 2014 $en_dash Present\",\n    \"PARTNER - Siege Analytics, Washington, DC | January 2014 $en_dash Present\"),\n  @(\"PRINCIPAL TECHNICAL ARCHITECT - Clarity and Rigour, Washington, DC | 2012 $en_dash 2014\",\n    \"DATA PRODUCTS MANAGER - Helm/Murmuration, Washington, DC | 2012 $en_dash 2014\"),\n  @(\"DIRECTOR OF TECHNOLOGY - Helm, Washington, DC | 2010 $en_dash 2012\",\n    \"SOFTWARE ENGINEER - Mautinoa Technologies, Washington, DC | 2010 $en_dash 2012\"),\n  @(\"SENIOR TECHNICAL ANALYST - GSD&M, Austin, TX | 2008 $en_dash 2010\",\n    \"SENIOR ANALYST - Myers Research, Washington, DC | 2008 $en_dash 2010\"),\n  @(\"TECHNICAL COORDINATOR - Progressive Change Campaign Committee, Washington, DC | 2006 $en_dash 2008\",\n    \"RESEARCH DIRECTOR - Progressive Change Campaign Committee, Washington, DC | 2006 $en_dash 2008\"),\n  @(\"TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 $en_dash 2004\",\n    \"INTERIM TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 $en_dash 2004\"),\n  @(\"TECHNICAL COORDINATOR - The Feldman Group, Washington, DC | 2000 $en_dash 2001\",\n    \"FIELD DIRECTOR - The Feldman Group, Washington, DC | 2000 $en_dash 2001\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.MatchSoundsLike = $false\n  $find.MatchAllWordForms = $false\n\n  $ok = $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n  if (-not $ok) {\n    throw \"Could not find text: $oldText\"\n  }\n}\n\nWrite-Output \"done\"\n"}
